$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3003.647
$ws.Range("I76").Value = 3003
$ws.Range("J76").Value = 3014
$ws.Range("K76").Value = 3003
$ws.Range("L76").Value = 3014
$ws.Range("M76").Value = -2688
$ws.Range("N76").Value = -3644

$ws.Range("H79").Value = 3003.647
$ws.Range("I79").Value = 3003
$ws.Range("J79").Value = 3014
$ws.Range("K79").Value = 3003
$ws.Range("L79").Value = 3014
$ws.Range("M79").Value = -1911
$ws.Range("N79").Value = -5198

$ws.Range("H80").Value = 47621748
$ws.Range("I80").Value = 66668250
$ws.Range("J80").Value = 5499.5
$ws.Range("K80").Value = 200004750
$ws.Range("L80").Value = 16498.5
$ws.Range("M80").Value = -200003752
$ws.Range("N80").Value = -18494.5

$ws.Range("H83").Value = 47621748
$ws.Range("I83").Value = 66668250
$ws.Range("J83").Value = 5499.5
$ws.Range("K83").Value = 600014250
$ws.Range("L83").Value = 49495.5
$ws.Range("M83").Value = -600009258
$ws.Range("N83").Value = -59479.5

$ws.Range("H86").Value = 6251870.5
$ws.Range("I86").Value = 9092138
$ws.Range("J86").Value = 3280.8
$ws.Range("K86").Value = 9092138
$ws.Range("L86").Value = 3280.8
$ws.Range("M86").Value = -9091015
$ws.Range("N86").Value = -5526.8

$ws.Range("H89").Value = 6251870.5
$ws.Range("I89").Value = 9092138
$ws.Range("J89").Value = 3280.8
$ws.Range("K89").Value = 45460690
$ws.Range("L89").Value = 16404
$ws.Range("M89").Value = -45455074
$ws.Range("N89").Value = -27636

$ws.Range("H137").Value = 2335348.5
$ws.Range("I137").Value = 5918026.5
$ws.Range("J137").Value = 6608.05
$ws.Range("K137").Value = 17754079.5
$ws.Range("L137").Value = 19824.15
$ws.Range("M137").Value = -17751529.5
$ws.Range("N137").Value = -24924.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3510.4666
$ws.Range("I61").Value = 2599
$ws.Range("J61").Value = 3650.6924
$ws.Range("K61").Value = 2599
$ws.Range("L61").Value = 3650.6924
$ws.Range("M61").Value = -2387
$ws.Range("N61").Value = -4074.6924

$ws.Range("H122").Value = 2368.4
$ws.Range("I122").Value = 2226.2856
$ws.Range("K122").Value = 6678.8568
$ws.Range("M122").Value = -4228.8568

$ws.Range("H133").Value = 26002.732
$ws.Range("J133").Value = 26002.732
$ws.Range("L133").Value = 26002.732
$ws.Range("N133").Value = -31062.732

$ws.Range("H135").Value = 22421.562
$ws.Range("J135").Value = 22421.562
$ws.Range("L135").Value = 22421.562
$ws.Range("N135").Value = -32561.562

$ws.Range("H136").Value = 3510.4666
$ws.Range("I136").Value = 2599
$ws.Range("J136").Value = 3650.6924
$ws.Range("K136").Value = 7797
$ws.Range("L136").Value = 10952.0772
$ws.Range("M136").Value = -5247
$ws.Range("N136").Value = -16052.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2350.65
$ws.Range("I94").Value = 2281.5334
$ws.Range("K94").Value = 2281.5334
$ws.Range("M94").Value = -1830.5334

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H134").Value = 3505.1592
$ws.Range("I134").Value = 3385
$ws.Range("J134").Value = 3561.2334
$ws.Range("K134").Value = 10155
$ws.Range("L134").Value = 10683.7002
$ws.Range("M134").Value = -7620
$ws.Range("N134").Value = -15753.7002

$ws.Range("H140").Value = 38103.8
$ws.Range("J140").Value = 38103.8
$ws.Range("L140").Value = 38103.8
$ws.Range("N140").Value = -48463.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9530518
$ws.Range("I31").Value = 2691.9092
$ws.Range("J31").Value = 13897438
$ws.Range("K31").Value = 2691.9092
$ws.Range("L31").Value = 13897438
$ws.Range("M31").Value = -2396.9092
$ws.Range("N31").Value = -13898028

$ws.Range("H34").Value = 9530518
$ws.Range("I34").Value = 2691.9092
$ws.Range("J34").Value = 13897438
$ws.Range("K34").Value = 2691.9092
$ws.Range("L34").Value = 13897438
$ws.Range("M34").Value = -2489.9092
$ws.Range("N34").Value = -13897842

$ws.Range("H86").Value = 3155.5881
$ws.Range("I86").Value = 3312.8
$ws.Range("J86").Value = 2931
$ws.Range("K86").Value = 3312.8
$ws.Range("L86").Value = 2931
$ws.Range("M86").Value = -2189.8
$ws.Range("N86").Value = -5177

$ws.Range("H89").Value = 3155.5881
$ws.Range("I89").Value = 3312.8
$ws.Range("J89").Value = 2931
$ws.Range("K89").Value = 16564
$ws.Range("L89").Value = 14655
$ws.Range("M89").Value = -10948
$ws.Range("N89").Value = -25887

$ws.Range("H93").Value = 17483.334
$ws.Range("I93").Value = 4950
$ws.Range("J93").Value = 19990
$ws.Range("K93").Value = 4950
$ws.Range("L93").Value = 19990
$ws.Range("M93").Value = -3078
$ws.Range("N93").Value = -23734

$ws.Range("H94").Value = 1068.3043
$ws.Range("J94").Value = 1273.8823
$ws.Range("L94").Value = 1273.8823
$ws.Range("N94").Value = -2175.8823

$ws.Range("H103").Value = 7472
$ws.Range("I103").Value = 7050.6665
$ws.Range("J103").Value = 10000
$ws.Range("K103").Value = 7050.6665
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = -5878.6665
$ws.Range("N103").Value = -12344

$ws.Range("H132").Value = 104851.57
$ws.Range("I132").Value = 2774.75
$ws.Range("J132").Value = 145682.3
$ws.Range("K132").Value = 8324.25
$ws.Range("L132").Value = 437046.9
$ws.Range("M132").Value = -5794.25
$ws.Range("N132").Value = -442106.9

$ws.Range("H134").Value = 621521.0600000001
$ws.Range("I134").Value = 742816.7
$ws.Range("J134").Value = 237418.33
$ws.Range("K134").Value = 2228450.1
$ws.Range("L134").Value = 712254.99
$ws.Range("M134").Value = -2225915.1
$ws.Range("N134").Value = -717324.99

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 597.86206
$ws.Range("I5").Value = 497.41666
$ws.Range("J5").Value = 1080
$ws.Range("K5").Value = 1492.24998
$ws.Range("L5").Value = 3240
$ws.Range("M5").Value = -1380.24998
$ws.Range("N5").Value = -3464

$ws.Range("H12").Value = 2604487.8
$ws.Range("I12").Value = 160.93333
$ws.Range("J12").Value = 4902423
$ws.Range("K12").Value = 482.79999
$ws.Range("L12").Value = 14707269
$ws.Range("M12").Value = -309.79999
$ws.Range("N12").Value = -14707615

$ws.Range("H113").Value = 855.46155
$ws.Range("I113").Value = 689.7273
$ws.Range("J113").Value = 977
$ws.Range("K113").Value = 2069.1819
$ws.Range("L113").Value = 2931
$ws.Range("M113").Value = 100.8181
$ws.Range("N113").Value = -7271

$ws.Range("H134").Value = 114447570
$ws.Range("I134").Value = 128752930
$ws.Range("J134").Value = 4700
$ws.Range("K134").Value = 386258790
$ws.Range("L134").Value = 14100
$ws.Range("M134").Value = -386253720
$ws.Range("N134").Value = -24240

$ws.Range("H135").Value = 597.86206
$ws.Range("I135").Value = 497.41666
$ws.Range("J135").Value = 1080
$ws.Range("K135").Value = 4476.74994
$ws.Range("L135").Value = 9720
$ws.Range("M135").Value = -1941.74994
$ws.Range("N135").Value = -14790

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1702.9
$ws.Range("I113").Value = 1685.5714
$ws.Range("J113").Value = 1743.3334
$ws.Range("K113").Value = 1685.5714
$ws.Range("L113").Value = 1743.3334
$ws.Range("M113").Value = 484.4286
$ws.Range("N113").Value = -6083.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2503.5
$ws.Range("J13").Value = 2503.5
$ws.Range("L13").Value = 2503.5
$ws.Range("N13").Value = -2783.5

$ws.Range("H22").Value = 1264
$ws.Range("I22").Value = 971.375
$ws.Range("J22").Value = 1654.1666
$ws.Range("K22").Value = 971.375
$ws.Range("L22").Value = 1654.1666
$ws.Range("M22").Value = -676.375
$ws.Range("N22").Value = -2244.1666

$ws.Range("H27").Value = 1264
$ws.Range("I27").Value = 971.375
$ws.Range("J27").Value = 1654.1666
$ws.Range("K27").Value = 971.375
$ws.Range("L27").Value = 1654.1666
$ws.Range("M27").Value = -864.375
$ws.Range("N27").Value = -1868.1666

$ws.Range("H40").Value = 2568.1875
$ws.Range("I40").Value = 2673.25
$ws.Range("J40").Value = 2253
$ws.Range("K40").Value = 2673.25
$ws.Range("L40").Value = 2253
$ws.Range("M40").Value = -2537.25
$ws.Range("N40").Value = -2525

$ws.Range("H55").Value = 557.54285
$ws.Range("I55").Value = 531.9091
$ws.Range("J55").Value = 600.9231
$ws.Range("K55").Value = 531.9091
$ws.Range("L55").Value = 600.9231
$ws.Range("M55").Value = -358.9091
$ws.Range("N55").Value = -946.9231

$ws.Range("H61").Value = 4280.6
$ws.Range("I61").Value = 4200.6665
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 4200.6665
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3998.6665
$ws.Range("N61").Value = -5404

$ws.Range("H68").Value = 3452.5386
$ws.Range("I68").Value = 2796
$ws.Range("K68").Value = 2796
$ws.Range("M68").Value = -2047

$ws.Range("H71").Value = 3452.5386
$ws.Range("I71").Value = 2796
$ws.Range("K71").Value = 13980
$ws.Range("M71").Value = -10236

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H113").Value = 4280.6
$ws.Range("I113").Value = 4200.6665
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4200.6665
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2030.6665
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 11112484
$ws.Range("I107").Value = 892
$ws.Range("K107").Value = 2676
$ws.Range("M107").Value = -756

$ws.Range("H132").Value = 2418198.8
$ws.Range("I132").Value = 4349683.5
$ws.Range("J132").Value = 3842.625
$ws.Range("K132").Value = 13049050.5
$ws.Range("L132").Value = 11527.875
$ws.Range("M132").Value = -13046520.5
$ws.Range("N132").Value = -16587.875

$ws.Range("H136").Value = 519444.6
$ws.Range("I136").Value = 707821.8
$ws.Range("K136").Value = 2123465.4
$ws.Range("M136").Value = -2120915.4
